# update to batpac v5
#
# Changes:
#  1. Move the "cathode | cathode binder | cathode binder (PVDF)" row from
#     row 27 down to just before the separator rows (new row 59).
#  2. Insert a new row "module | module packaging | module polymer panels"
#     right after the "module container" row, with a taller row height and
#     a Segoe UI / black font applied to the new component-name cell.
#  3. Rename "module compression plates" to "module row rack".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Relocate the cathode binder (PVDF) row ---------------------------
# It currently lives at row 27 (A27:C27 = cathode / cathode binder / cathode
# binder (PVDF)). Remove it from there (rows below shift up), then insert a
# fresh row right before the "coated separator (5um+2um)" row and fill it
# back in with the same three values.

$ws.Rows.Item(27).Delete()

# After the delete, the separator block ("coated separator (5um+2um)" etc.)
# starts at row 58. Insert a blank row there to make room, which pushes the
# separator block down to start at row 59, and put the cathode binder row
# in the freed row 58.
$ws.Rows.Item(58).Insert()
$ws.Cells.Item(58, 1).Value2 = "cathode"
$ws.Cells.Item(58, 2).Value2 = "cathode binder"
$ws.Cells.Item(58, 3).Value2 = "cathode binder (PVDF)"

# --- 2. Insert the new "module polymer panels" row ------------------------
# "module container" is now at row 48 (module / module packaging / module
# container). Insert a new row right after it.
$ws.Rows.Item(49).Insert()
$ws.Cells.Item(49, 1).Value2 = "module"
$ws.Cells.Item(49, 2).Value2 = "module packaging"
$ws.Cells.Item(49, 3).Value2 = "module polymer panels"

# New row gets a taller height and a distinct font on the component cell.
$ws.Rows.Item(49).RowHeight = 16.5
$f = $ws.Cells.Item(49, 3).Font
$f.Name = "Segoe UI"
$f.Color = 0

# --- 3. Rename "module compression plates" -> "module row rack" ----------
$ws.Cells.Item(54, 3).Value2 = "module row rack"

# --- Cosmetic sheet-level touch-ups matching the saved workbook ----------
$ws.Range("I18").Select()
